# Generate Report for Handoff
# Replaces the localized file's UUID-based name (8afccafa-...  -> 9fa47fb9-...)
# throughout the handoff status workbook, and refreshes the associated
# handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldGuid = "8afccafa-1d18-48ae-a109-41c76d492241"
$newGuid = "9fa47fb9-3243-4839-8df9-1b5505b5a8bf"

$newMdName      = "$newGuid.md"
$newMdDisplay   = "e2e\$newGuid.md"

$newZhXlfName = "$newGuid.182727f9ca46ddb1006d71db2ac20dc2beb218cb.zh-cn.xlf"
$newDeXlfName = "$newGuid.182727f9ca46ddb1006d71db2ac20dc2beb218cb.de-de.xlf"

$newHoDate      = "2016-09-08 05:16:19"
$newZhHoDate    = "2016-09-08 05:16:14"

$oldHyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dfcb235a9bc483140808c0e3971297e7553bf8ce/e2e/$oldGuid.md"

# ---- Overview sheet ----
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdDisplay
$wsOverview.Range("G2").Value = $newHoDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $oldHyperlinkAddress, "", "", $newMdDisplay)

# ---- zh-cn sheet ----
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("G2").Value = $newZhXlfName
$wsZhCn.Range("H2").Value = $newZhHoDate
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $oldHyperlinkAddress, "", "", $newMdName)

# ---- de-de sheet ----
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("G2").Value = $newDeXlfName
$wsDeDe.Range("H2").Value = $newHoDate
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $oldHyperlinkAddress, "", "", $newMdName)
